# "updated docs and default settings in line"
# Update the default parameter values on row 2 of "sheet 1" and move the
# active selection from J3 to M3 (one column to the right of the last
# parameter column), matching the author's in-line edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet 1")

# training_block_length_multiplier: 1 -> 5
$ws.Range("G2").Value = 5

# training_criterion: 8 -> 38
$ws.Range("I2").Value = 38

# testing_criterion: 8 -> 7
$ws.Range("J2").Value = 7

# max_training: 10 -> 6
$ws.Range("K2").Value = 6

# max_training_and_testing: 2 -> 3
$ws.Range("M2").Value = 3

# Move the active cell/selection to M3 (was J3)
$ws.Range("M3").Select()
